$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = $false
$ws.Range("B1").Value = "simple"
$ws.Range("C1").Value = 58.25483664188928

$ws.Range("A2").Value = $true
$ws.Range("B2").Value = "simple"
$ws.Range("C2").Value = 51.23291552776303

$ws.Range("A3").Value = $false
$ws.Range("B3").Value = "hyper_heuristic"
$ws.Range("C3").Value = 104.828563401099

$ws.Range("A4").Value = $true
$ws.Range("B4").Value = "hyper_heuristic"
$ws.Range("C4").Value = 78.52288603930523
